$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-apply the "Text" number format to the header row and the row-label
# column, matching the regenerated style table from the source CSV export.
$ws.Range("A1:C1").NumberFormat = "@"
$ws.Range("A2:A11").NumberFormat = "@"

# Updated prediction-score values (column B) produced by the refreshed
# quadratic-svm run recorded in ful-path.csv.
$ws.Range("B2").Value = 0.018851785135483468
$ws.Range("B3").Value = -0.18916947023327602
$ws.Range("B4").Value = -0.28751753218624998
$ws.Range("B5").Value = -0.15706924481213314
$ws.Range("B6").Value = -0.0799342868630184
$ws.Range("B7").Value = -0.1850295281991059
$ws.Range("B8").Value = -0.23046933680048198
$ws.Range("B9").Value = -0.30976566846797837
$ws.Range("B10").Value = -0.14530143628924108
$ws.Range("B11").Value = 0.042806602233532942
